$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '20.560.17'
$ws.Range("E2").Value = '  -0.06%  '

# Row 3
$ws.Range("D3").Value = '1.481.07'
$ws.Range("E3").Value = '  +0.58%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9768'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.88%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '279.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.80%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3663'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.67%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3079'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.90%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.64'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.11%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.067'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06668'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.90%  '

# Row 12
$ws.Range("E12").Value = '  -0.18%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.528'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.76%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.213'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.64%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9768'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.60%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001034'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.62%  '

# Row 18
$ws.Range("D18").Value = '1.479.21'
$ws.Range("E18").Value = '  +0.00%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05937'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.505'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.85%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.27%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.93%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.249'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.13%  '

# Row 25
$ws.Range("D25").Value = '20.616.99'
$ws.Range("E25").Value = '  -0.24%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.45%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.141'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.67%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.07%  '

# Row 29
$ws.Range("D29").Value = '1.637.18'
$ws.Range("E29").Value = '  -0.39%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.974'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8205'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.43%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.42%  '

# Row 34
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08039'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.14%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.546'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.57%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.222'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.17%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05863'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.67%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.731'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.51%  '

# Row 39
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.846'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.42%  '

# Row 40
$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9762'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.60%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02050'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.66%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.02%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1890'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5308'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.96%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.529'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.00%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.64%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.68%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5206'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.34%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.805'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.83%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06476'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.37%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9963'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.69%  '
